# Applies updated odds values per row (3, 8, 11, 23, 27) as per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 1.83
$ws.Range("H3").Value = 3.3
$ws.Range("I3").Value = 4.33
$ws.Range("X3").Value = 8
$ws.Range("Z3").Value = 15
$ws.Range("AC3").Value = 8
$ws.Range("AH3").Value = 21
$ws.Range("AI3").Value = 15
$ws.Range("AL3").Value = 41
$ws.Range("AN3").Value = 3.75
$ws.Range("AR3").Value = 67
$ws.Range("AW3").Value = 6
$ws.Range("AY3").Value = 34
$ws.Range("BA3").Value = 126
$ws.Range("BB3").Value = 301
# Row 8
$ws.Range("G8").Value = 4.15
$ws.Range("H8").Value = 3.6
$ws.Range("I8").Value = 1.75
$ws.Range("J8").Value = 4.45
$ws.Range("K8").Value = 2.18
$ws.Range("L8").Value = 2.3
$ws.Range("M8").Value = 1.03
$ws.Range("N8").Value = 11
$ws.Range("Q8").Value = 1.7
$ws.Range("R8").Value = 1.91
$ws.Range("U8").Value = 1.65
$ws.Range("V8").Value = 1.98
$ws.Range("W8").Value = 13
$ws.Range("X8").Value = 25
$ws.Range("Y8").Value = 13.5
$ws.Range("Z8").Value = 65
$ws.Range("AA8").Value = 37
$ws.Range("AB8").Value = 40
$ws.Range("AC8").Value = 11.5
$ws.Range("AD8").Value = 7.1
$ws.Range("AE8").Value = 14
$ws.Range("AF8").Value = 60
$ws.Range("AG8").Value = 7.7
$ws.Range("AH8").Value = 8.75
$ws.Range("AI8").Value = 8
$ws.Range("AJ8").Value = 14.5
$ws.Range("AK8").Value = 13.5
$ws.Range("AN8").Value = 6
$ws.Range("AO8").Value = 23
$ws.Range("AP8").Value = 28
$ws.Range("AQ8").Value = 120
$ws.Range("AR8").Value = 150
$ws.Range("AS8").Value = 350
$ws.Range("AT8").Value = 2.65
$ws.Range("AU8").Value = 7.1
$ws.Range("AW8").Value = 3.65
$ws.Range("AX8").Value = 8.5
$ws.Range("AY8").Value = 17
$ws.Range("AZ8").Value = 29
$ws.Range("BA8").Value = 55
# Row 11
$ws.Range("O11").Value = 1.29
$ws.Range("P11").Value = 3.5
$ws.Range("Q11").Value = 1.95
$ws.Range("R11").Value = 1.85
# Row 23
$ws.Range("G23").Value = 3.25
$ws.Range("I23").Value = 2.1
$ws.Range("J23").Value = 3.5
$ws.Range("L23").Value = 2.6
$ws.Range("M23").Value = 1.02
$ws.Range("N23").Value = 19
$ws.Range("S23").Value = 1.25
$ws.Range("T23").Value = 3.75
$ws.Range("U23").Value = 1.44
$ws.Range("V23").Value = 2.63
$ws.Range("X23").Value = 21
$ws.Range("Y23").Value = 12
$ws.Range("AC23").Value = 19
$ws.Range("AL23").Value = 19
$ws.Range("AO23").Value = 17
$ws.Range("AP23").Value = 21
$ws.Range("AQ23").Value = 51
$ws.Range("AT23").Value = 3.75
$ws.Range("BB23").Value = 81
# Row 27
$ws.Range("G27").Value = 2.2
$ws.Range("H27").Value = 2.95
$ws.Range("I27").Value = 3.35
$ws.Range("J27").Value = 2.85
$ws.Range("K27").Value = 1.95
$ws.Range("L27").Value = 3.9
$ws.Range("O27").Value = 1.45
$ws.Range("P27").Value = 2.37
$ws.Range("Q27").Value = 2.32
$ws.Range("S27").Value = 1.47
$ws.Range("T27").Value = 2.32
$ws.Range("U27").Value = 1.98
$ws.Range("V27").Value = 1.65
$ws.Range("W27").Value = 5.9
$ws.Range("X27").Value = 9.25
$ws.Range("Y27").Value = 9.5
$ws.Range("Z27").Value = 21
$ws.Range("AB27").Value = 40
$ws.Range("AC27").Value = 6.8
$ws.Range("AD27").Value = 5.8
$ws.Range("AE27").Value = 17
$ws.Range("AF27").Value = 100
$ws.Range("AG27").Value = 8
$ws.Range("AH27").Value = 16.5
$ws.Range("AJ27").Value = 50
$ws.Range("AK27").Value = 35
$ws.Range("AN27").Value = 3.9
$ws.Range("AO27").Value = 11.75
$ws.Range("AP27").Value = 23
$ws.Range("AQ27").Value = 50
$ws.Range("AT27").Value = 2.3
$ws.Range("AU27").Value = 7.4
$ws.Range("AV27").Value = 75
$ws.Range("AW27").Value = 5
$ws.Range("AX27").Value = 19
$ws.Range("AZ27").Value = 100

Write-Host "Applied 114 cell updates"
